# Auto-generated edit script: update crypto price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; D="69.296.85"; E="  -0.36%  "},
    @{Row=3; D="2.489.13"; E="  -0.97%  "},
    @{Row=4; E="  -0.03%  "},
    @{Row=5; D="567.97"; E="  -0.83%  "},
    @{Row=6; D="165.60"; E="  -0.29%  "},
    @{Row=8; D="0.510"; E="  -0.64%  "},
    @{Row=9; D="0.158"; E="  -0.26%  "},
    @{Row=10; E="  -0.90%  "},
    @{Row=11; D="0.346"; E="  -2.96%  "},
    @{Row=12; D="4.87"; E="  -0.62%  "},
    @{Row=13; D="2.944.20"; E="  -1.03%  "},
    @{Row=14; D="69.203.08"; E="  -0.33%  "},
    @{Row=15; D="0.0000174"; E="  -1.06%  "},
    @{Row=16; D="24.08"; E="  -2.93%  "},
    @{Row=17; D="2.463.23"; E="  -2.56%  "},
    @{Row=18; D="11.16"; E="  -1.06%  "},
    @{Row=19; B="BitcoinCash"; C="https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"; D="353.14"; E="  +1.03%  "},
    @{Row=20; B="Uniswap"; C="https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"; D="7.34"; E="  -3.30%  "},
    @{Row=21; D="3.89"; E="  -0.23%  "},
    @{Row=22; D="1.92"; E="  -3.14%  "},
    @{Row=23; E="  -0.02%  "},
    @{Row=24; D="69.19"; E="  -1.64%  "},
    @{Row=25; D="3.78"; E="  -3.28%  "},
    @{Row=26; E="  -0.81%  "},
    @{Row=27; D="8.62"; E="  -3.09%  "},
    @{Row=28; D="0.996"; E="  -0.37%  "},
    @{Row=29; D="0.0₃0867"; E="  -2.28%  "},
    @{Row=30; D="7.49"; E="  -3.87%  "},
    @{Row=31; D="3.54"; E="  +137.26%  "},
    @{Row=32; B="Bittensor"; C="https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"; D="437.89"; E="  -5.27%  "},
    @{Row=33; B="Fetch.AI"; C="https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"; D="1.19"; E="  -3.53%  "},
    @{Row=34; D="0.999"; E="  -0.06%  "},
    @{Row=35; D="1.70"; E="  -1.25%  "},
    @{Row=36; B="Monero"; C="https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; D="154.13"; E="  -2.02%  "},
    @{Row=37; B="Kaspa"; C="https://coinranking.com/coin/V8GxkwWow+kaspa-kas"; D="0.112"; E="  -3.52%  "},
    @{Row=38; D="19.05"; E="  -0.22%  "},
    @{Row=39; D="18.09"; E="  -1.99%  "},
    @{Row=40; E="  -0.04%  "},
    @{Row=41; D="0.312"; E="  -1.78%  "},
    @{Row=42; D="4.57"; E="  -2.44%  "},
    @{Row=43; D="1.56"; E="  -2.21%  "},
    @{Row=44; E="  -1.99%  "},
    @{Row=45; E="  -4.09%  "},
    @{Row=46; D="138.64"; E="  -2.25%  "},
    @{Row=47; D="3.42"; E="  -1.20%  "},
    @{Row=48; D="0.503"; E="  -3.08%  "},
    @{Row=49; E="  -1.22%  "},
    @{Row=50; D="0.572"; E="  -0.97%  "},
    @{Row=51; E="  -0.46%  "},
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey("B")) { $ws.Range("B" + $r).Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Range("C" + $r).Value = $u.C }
    if ($u.ContainsKey("D")) {
        $ws.Range("D" + $r).NumberFormat = "@"
        $ws.Range("D" + $r).Value = $u.D
    }
    if ($u.ContainsKey("E")) { $ws.Range("E" + $r).Value = $u.E }
}
